$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The betting-odds columns (F:V) for several row pairs were re-scraped and
# ended up swapped between two rows that share the same match date. Restore
# the correct pairing by swapping columns F:V between each pair of rows.
$swapPairs = @(
    @(34, 36),
    @(35, 37),
    @(39, 40),
    @(41, 42),
    @(56, 57)
)

foreach ($pair in $swapPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]
    $rng1 = $ws.Range("F$r1`:V$r1")
    $rng2 = $ws.Range("F$r2`:V$r2")
    $v1 = $rng1.Value()
    $v2 = $rng2.Value()
    $rng1.Value = $v2
    $rng2.Value = $v1
}

# Append a new match row (146) at the end of the table, copying the
# formatting (styles) of the last existing row (145) first.
$ws.Range("A145:V145").Copy($ws.Range("A146:V146"))

$ws.Range("A146").Value = 145
$ws.Range("E146").Value = 45241.67708333334
$ws.Range("F146").Value = "Olimpia Elblag"
$ws.Range("G146").Value = 2
$ws.Range("H146").Value = "Chojniczanka"
$ws.Range("I146").Value = 1
$ws.Range("J146").Value = 2.31
$ws.Range("K146").Value = "10/11/2023 04:42"
$ws.Range("L146").Value = 2.22
$ws.Range("M146").Value = "11/11/2023 15:59"
$ws.Range("N146").Value = 3.02
$ws.Range("O146").Value = "10/11/2023 04:42"
$ws.Range("P146").Value = 3.23
$ws.Range("Q146").Value = "11/11/2023 15:59"
$ws.Range("R146").Value = 2.85
$ws.Range("S146").Value = "10/11/2023 04:42"
$ws.Range("T146").Value = 3.21
$ws.Range("U146").Value = "11/11/2023 15:59"
$ws.Range("V146").Value = "https://www.betexplorer.com/football/poland/division-2/olimpia-elblag-chojniczanka/EexfDkKk/"
